$wb = $excel.ActiveWorkbook

# The new weekly sheet has the same layout as the previous week's sheet, so
# duplicate it (keeps header style, outline/page setup props, column widths, etc.)
# and then overwrite its ranking data in place.
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$lastSheet.Copy([System.Reflection.Missing]::Value, $lastSheet)

$ws = $wb.Worksheets.Item($sheetCount + 1)
$ws.Name = "magapoke_2025-12-03"

# Header row (unchanged: rank / title)
$ws.Range("A1").Value = "rank"
$ws.Range("B1").Value = "title"

# This week's ranking data: title ordered by rank 1..45
$titles = @(
    "スルガメテオ",
    "K-9~警視庁公安部公安第9課異能対策係~",
    "ドリーム☆ジャンボ☆ガール",
    "アイドラトリィ",
    "黄昏町プリズナーズ",
    "ナキナギ",
    "生きたがりの人狼",
    "ゼロとヒャク",
    "黒月のイェルクナハト",
    "普通の本はありません！",
    "その青春",
    "夜鐘のキト",
    "お母さん冒険者、ログインボーナスでスキル【主婦】に目覚めました。週一貰えるチラシで冒険者生活頑張ります！",
    "ハードワーカー中田",
    "春くらり",
    "屋根の下のアルテミス",
    "君が監督！",
    "卒業アルバムの彼女たち",
    "MYS",
    "篝家の８兄弟",
    "限界集落を脱村した錬金術士、都会で""最強""なのがバレまくる。～老害どもにはいい加減愛想が尽きました～",
    "異世界グルメで成り上がり無双～山に追放されたので、のんびりキャンプを楽しんでいたらいつの間にか強くなっていて、王侯貴族や実力者たちが俺を放っておいてくれません。一方、俺を追放した貴族たちは破滅が始まる～",
    "それがメイドのカンナです",
    "せいぶつ部の田辺くん",
    "ハナバス　苔石花江のバスケ論",
    "明智ナンバーワン",
    "ルックスＹを選んでしまいました 〜やり込んでいるゲームに転生したはずなのに、未実装のガチャで攻略をすることになった件〜",
    "ナマイキ旭ちゃんをわからせたい",
    "ともだちづくり",
    "じゅーくぼっくす",
    "追放されなかった男　～二度目の人生は土下座から始まりました～",
    "鳴るさんだぁ",
    "平成転生",
    "永久のユウグレ",
    "皇女転生　～伝説の大魔導士（♂）、姫騎士となりて伝説の令嬢騎士団を作り無双する～",
    "JK Biker",
    "白銀のキュイジーヌ～明治外交官の料理人～",
    "人生逆転ダンジョン",
    "ハプスブルク家の華麗なる受難",
    "鉱石令嬢〜没落した悪役令嬢が炭鉱で一山当てるまでのお話〜",
    "眠れる森のレガ",
    "花子狩り",
    "〈小市民〉 春期限定いちごタルト事件",
    "東京デスレース",
    "イエティ、とある日々"
)

for ($i = 0; $i -lt $titles.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i + 1
    $ws.Cells.Item($row, 2).Value = $titles[$i]
}

